# Actualizar fixtures y agregar nuevas ligas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New fixture rows to append starting at row 122
$rows = @(
    @{ A="2025-07-12"; B="Curico Unido";              C="San Marcos de Arica";       D=5;  E=3; F=1348357; G=2;  H=4; I=1; J=3; K=0; L=0; M=1; N=2; O=4; P=1; Q=54; R=46; S="L" },
    @{ A="2025-07-12"; B="Union San Felipe";           C="Cobreloa";                  D=3;  E=0; F=1348354; G=3;  H=5; I=2; J=1; K=0; L=0; M=1; N=0; O=2; P=0; Q=48; R=52; S="L" },
    @{ A="2025-07-12"; B="Rangers de Talca";           C="Deportes Santa Cruz";       D=2;  E=2; F=1348358; G=1;  H=5; I=4; J=3; K=1; L=0; M=1; N=1; O=1; P=1; Q=50; R=50; S="E" },
    @{ A="2025-07-12"; B="Universidad de Concepcion";  C="Santiago Wanderers";        D=0;  E=3; F=1348359; G=10; H=1; I=4; J=2; K=0; L=1; M=0; N=2; O=0; P=1; Q=47; R=53; S="V" },
    @{ A="2025-07-13"; B="Antofagasta";                C="San Luis";                  D=1;  E=2; F=1348353; G=12; H=3; I=3; J=4; K=1; L=0; M=1; N=1; O=0; P=1; Q=53; R=47; S="V" },
    @{ A="2025-07-13"; B="Santiago Morning";           C="Deportes Copiapo";          D=1;  E=0; F=1348356; G=2;  H=5; I=4; J=2; K=0; L=0; M=0; N=0; O=1; P=0; Q=43; R=57; S="L" },
    @{ A="2025-07-13"; B="Recoleta";                   C="Magallanes";                D=2;  E=0; F=1348355; G=3;  H=6; I=2; J=3; K=0; L=1; M=1; N=0; O=1; P=0; Q=46; R=54; S="L" },
    @{ A="2025-07-13"; B="Deportes Temuco";            C="Concepción";                D=2;  E=1; F=1348360; G=3;  H=8; I=3; J=2; K=0; L=2; M=1; N=0; O=1; P=1; Q=53; R=47; S="L" }
)

$startRow = 122
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S")

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    foreach ($col in $cols) {
        $cell = $ws.Range("$col$r")
        if ($col -eq "A") {
            # Column A holds date-like text (e.g. "2025-07-12"); force it to stay
            # plain text instead of being auto-converted into a date serial number.
            $cell.NumberFormat = "@"
            $cell.Value = $row[$col]
            $cell.ClearFormats()
        } else {
            $cell.Value = $row[$col]
        }
    }
}

# Update view: scroll position and selection on the newly added data
$ws.Application.ActiveWindow.ScrollRow = 115
$ws.Range("O129").Select()

# Column A no longer needs a custom width; reset it to the sheet's default width
$ws.Columns.Item(1).ColumnWidth = $ws.Columns.Item(4).ColumnWidth
